$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptocurrency price and 1h-volume-change data for each row.
# Map of row -> (Price, Volume(1h)) values, applied as literal text
# (matching the source data which stores these as text, not numbers).
$updates = @{
    2 = @("273.48", "-1.77%")
    3 = @("26.76", "-2.19%")
    4 = @("4.904", "2.25%")
    5 = @("0.06333", "1.14%")
    6 = @("6.895", "0.81%")
    7 = @("3.362", "3.39%")
    8 = @("1.281", "35.78%")
    9 = @("0.8796", "0.25%")
    10 = @("0.1456", "-0.42%")
    11 = @("0.04963", "-4.18%")
    12 = @("0.07396", "1.02%")
    13 = @("0.03121", "-0.98%")
    14 = @("0.09058", "0.10%")
    15 = @("0.001595", "2.67%")
    16 = @("0.0006334", "1.25%")
    17 = @("0.006025", "4.62%")
    18 = @("3.466", "-0.32%")
    19 = @("2.271", "-0.16%")
    21 = @("0.1328", "1.61%")
    22 = @("3.918", "2.08%")
    23 = @("0.04415", "2.12%")
    24 = @("0.001180", "0.44%")
    25 = @("0.003697", "-13.17%")
    26 = @("0.0001203", "0.41%")
    27 = @("0.0001705", "1.19%")
    40 = @("0.04053", "0.32%")
    41 = @("0.006662", "6.50%")
    42 = @("0.1168", "1.58%")
    43 = @("0.002106", "-0.84%")
    44 = @("0.01198", "-11.79%")
    45 = @("0.00005337", "3.96%")
    46 = @("2.356", "22.85%")
    47 = @("0.02005", "-32.87%")
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $dCell = $ws.Range("D$row")
    $eCell = $ws.Range("E$row")

    # Force text number-format so Excel stores the values as literal
    # strings (preserving formatting like trailing zeros and "%" signs)
    # instead of auto-converting them to numeric cell values.
    $dCell.NumberFormat = "@"
    $dCell.Value = $vals[0]
    $dCell.Style = "Normal"

    $eCell.NumberFormat = "@"
    $eCell.Value = $vals[1]
    $eCell.Style = "Normal"
}
